# The deck originally had 3 slides:
#   1) sldId 256 - "Group 133" word-chain diagram (rear/odec/deco/area/code/dear)
#   2) sldId 261 - an exact duplicate of slide 1's "Group 133" diagram
#   3) sldId 260 - the TSP / orca-car-arc slide
#
# The commit removes the duplicate second slide, leaving just the
# original diagram slide followed by the TSP slide.

$p = $ppt.ActivePresentation

$s = $p.Slides.Item(2)
$s.Delete()
